$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet currently holds one tab-separated string per row in column A
# ("code\tsubject\tbody", "st\tgo-sns started\t...", ...). Split each row
# into three real columns (code / subject / body).
# ---------------------------------------------------------------------------
$data = @(
  @("code", "subject", "body"),
  @("st", "go-sns started", "go-sns has started successfully"),
  @("sh", "A system is shutting down", "The attached system is attempting to shutdown"),
  @("tc", "Task Complete", "A task has been completed"),
  @("tf", "Task Failed", "A task has failed")
)

for ($r = 0; $r -lt $data.Length; $r++) {
  for ($c = 0; $c -lt 3; $c++) {
    $ws.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
  }
}

# ---------------------------------------------------------------------------
# Apply the new font (Helvetica Neue, 10pt, black) to the whole table.
# Build it once on a scratch cell so the font/style only needs a single
# resolved entry, then fan it out with a format-only paste.
# ---------------------------------------------------------------------------
$scratch = $ws.Range("Z100")
$scratch.Font.Color = 0
$scratch.Font.Size = 10
$scratch.Font.Name = "Helvetica Neue"
$scratch.Copy()

$table = $ws.Range("A1:C5")
$table.PasteSpecial(-4122)
$scratch.Clear()
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Column B needs to be wide enough to show the longest subject line.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 20.6

# ---------------------------------------------------------------------------
# Restore the selection to where the author left it.
# ---------------------------------------------------------------------------
$ws.Range("F12").Select()
